$d = $word.ActiveDocument

# --- 1) Strike-through the "Synkroniser repositorys" paragraph -------------
$pSync = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Synkroniser*repositorys*") {
        $pSync = $p
        break
    }
}
if ($pSync -ne $null) {
    $pSync.Range.Font.StrikeThrough = 1
}

# --- 2) Strike-through the "Opdater burndown-chart" paragraph --------------
#        and drop the _GoBack bookmark that currently lives on it
#        (it will be re-created further below, on the "Faerdiggoere..." line,
#        which also removes it from here since _GoBack is unique).
$pBurn = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Opdater burndown-chart*") {
        $pBurn = $p
        break
    }
}
if ($pBurn -ne $null) {
    $pBurn.Range.Font.StrikeThrough = 1
}

# --- 3) "Faerdiggoere afsnit om sprint review og retrospekt (casper)." -----
#        -> remove " (casper)" and move the _GoBack bookmark to sit right
#        before the trailing ".".
$pDone = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*sprint review og retrospekt*") {
        $pDone = $p
        break
    }
}
if ($pDone -ne $null) {
    $pStart = $pDone.Range.Start
    $txt = $pDone.Range.Text
    $removeStart = $txt.IndexOf(" (casper)")
    if ($removeStart -ge 0) {
        $delRange = $d.Range($pStart + $removeStart, $pStart + $removeStart + " (casper)".Length)
        $delRange.Delete()
    }

    $pDone2 = $pDone
    $pStart2 = $pDone2.Range.Start
    $txt2 = $pDone2.Range.Text
    $dotOffset = $txt2.LastIndexOf(".")
    $bmRange = $d.Range($pStart2 + $dotOffset, $pStart2 + $dotOffset)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
